# This script reproduces the Sat Apr 8 18:46:24 UTC 2023 GitHub Actions refresh of the
# "cryptos" price table: coin rankings 11-25/28-29 shifted by one position (OKB dropped,
# LEO added) and every row's Price (D) / Volume(1h) (E) value was refreshed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-51. UpdateBC marks rows whose Coin/Link columns actually
# changed (the rest keep their existing Coin/Link and only get new Price/Volume).
$rows = @(
    [PSCustomObject]@{ Row = 2; B = "Bitcoin"; C = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D = "28.558.73"; E = "  +2.00%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 3; B = "Ethereum"; C = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D = "1.894.53"; E = "  +1.86%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 4; B = "TetherUSD"; C = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D = "1.038"; E = "  +3.37%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 5; B = "BNB"; C = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D = "318.46"; E = "  +1.93%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 6; B = "USDC"; C = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D = "1.021"; E = "  +1.79%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 7; B = "XRP"; C = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D = "0.5171"; E = "  +0.64%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 8; B = "Cardano"; C = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D = "0.3947"; E = "  +2.98%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 9; B = "Dogecoin"; C = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D = "0.08399"; E = "  +1.69%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 10; B = "Polygon"; C = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D = "1.130"; E = "  +1.81%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 11; B = "Polkadot"; C = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D = "6.279"; E = "  +1.36%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 12; B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "1.902.28"; E = "  +1.91%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 13; B = "Solana"; C = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D = "20.57"; E = "  -0.09%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 14; B = "Chainlink"; C = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D = "7.305"; E = "  +0.78%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 15; B = "BinanceUSD"; C = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D = "1.024"; E = "  +2.04%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 16; B = "ShibaInu"; C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D = "0.00001112"; E = "  +1.31%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 17; B = "Litecoin"; C = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D = "91.25"; E = "  +0.74%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 18; B = "TRON"; C = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D = "0.06756"; E = "  +1.51%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 19; B = "Avalanche"; C = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D = "17.93"; E = "  +1.38%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 20; B = "Dai"; C = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; D = "1.021"; E = "  +1.84%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 21; B = "Uniswap"; C = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D = "6.051"; E = "  +0.72%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 22; B = "WrappedBTC"; C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D = "28.505.13"; E = "  +1.69%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 23; B = "Cosmos"; C = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D = "11.24"; E = "  +1.28%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 24; B = "Toncoin"; C = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D = "2.305"; E = "  +2.26%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 25; B = "LEO"; C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D = "3.410"; E = "  +0.61%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 26; B = "WrappedliquidstakedEther2.0"; C = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D = "2.124.45"; E = "  +2.29%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 27; B = "Monero"; C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D = "161.78"; E = "  +2.33%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 28; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "20.93"; E = "  +2.29%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 29; B = "LidoDAOToken"; C = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D = "2.460"; E = "  -2.27%  "; UpdateBC = $true },
    [PSCustomObject]@{ Row = 30; B = "BitcoinCash"; C = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D = "127.07"; E = "  +1.93%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 31; B = "Stellar"; C = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D = "0.1065"; E = "  +0.08%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 32; B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "1.046"; E = "  +1.76%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 33; B = "Filecoin"; C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D = "5.954"; E = "  -0.46%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 34; B = "HuobiToken"; C = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D = "3.659"; E = "  +1.59%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 35; B = "FraxShare"; C = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D = "9.467"; E = "  +1.23%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 36; B = "VeChain"; C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D = "0.02461"; E = "  +1.89%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 37; B = "Hedera"; C = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; D = "0.06625"; E = "  +2.03%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 38; B = "Algorand"; C = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D = "0.2228"; E = "  +2.74%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 39; B = "TheSandbox"; C = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D = "0.6521"; E = "  -0.51%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 40; B = "TrustWalletToken"; C = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D = "1.256"; E = "  +2.48%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 41; B = "ARBITRUM"; C = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D = "1.191"; E = "  -0.38%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 42; B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "5.040"; E = "  +0.72%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 43; B = "Aptos"; C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D = "11.23"; E = "  +0.73%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 44; B = "Decentraland"; C = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; D = "0.6152"; E = "  +0.04%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 45; B = "EnergySwap"; C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D = "13.23"; E = "  +1.39%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 46; B = "PancakeSwap"; C = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D = "3.734"; E = "  +2.05%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 47; B = "WEMIXTOKEN"; C = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; D = "1.295"; E = "  +1.14%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 48; B = "EOS"; C = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"; D = "1.243"; E = "  +2.29%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 49; B = "NEARProtocol"; C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D = "2.019"; E = "  +0.57%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 50; B = "Quant"; C = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D = "122.18"; E = "  +1.45%  "; UpdateBC = $false },
    [PSCustomObject]@{ Row = 51; B = "Cronos"; C = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D = "0.06940"; E = "  +1.43%  "; UpdateBC = $false }
)

# Column D holds price text that can look numeric (e.g. "1.038" or "3.410"); force
# text format on the whole price column range first so Excel does not silently
# convert it to a number/date when the Value is assigned below.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

foreach ($r in $rows) {
    if ($r.UpdateBC) {
        $ws.Cells.Item($r.Row, 2).Value = $r.B
        $ws.Cells.Item($r.Row, 3).Value = $r.C
    }
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

# Restore the original (General) number format/style on the price column so the
# workbook style table matches the source layout instead of keeping the temporary
# text format used above.
$priceRange.Style = $ws.Range("A1").Style

